$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.659.36"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.602.80"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'212.19"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'27.10"
$ws.Range("E8").Value = "  +9.32%  "
$ws.Range("D9").Value = "'43.38"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "1.832.00"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").Value = "1.613.47"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").Value = "29.647.62"
$ws.Range("E15").Value = "  +3.43%  "
$ws.Range("D16").Value = "'0.535"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "'63.29"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "'241.00"
$ws.Range("E19").Value = "  +6.13%  "
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("D21").Value = "0.0₃0691"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'3.99"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'2.07"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").Value = "'154.62"
$ws.Range("D27").Value = "'15.33"
$ws.Range("E27").Value = "  +3.65%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D35").Value = "1.428.71"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "'1.53"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").Value = "'0.536"
$ws.Range("E41").Value = "  +3.59%  "
$ws.Range("B42").Value = "BitcoinSV"
$ws.Range("C42").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D42").Value = "'54.78"
$ws.Range("E42").Value = "  +29.97%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.97"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "'0.0487"
$ws.Range("E44").Value = "  +6.19%  "
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'65.86"
$ws.Range("D48").Value = "'0.950"
$ws.Range("E48").Value = "  +13.12%  "
$ws.Range("D49").Value = "'5.27"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("D50").Value = "1.743.20"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'86.42"
$ws.Range("E51").Value = "  +2.06%  "
